# Update "Return_with_prediction" (G), "return_pct_change" (H), and
# "mean_return_pct_change" (I, row 2 only) with newly predicted values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.10559013179121
$ws.Range("H2").Value = -19.65489589087305
$ws.Range("I2").Value = 3.117300874808048

$ws.Range("G3").Value = 0.1058412693455825
$ws.Range("H3").Value = 18.87309321146904

$ws.Range("G4").Value = -0.5537239187016033
$ws.Range("H4").Value = 10.01384164028147

$ws.Range("G5").Value = -0.6266074489483778
$ws.Range("H5").Value = -2.622728685391902

$ws.Range("G6").Value = 0.2206606219079868
$ws.Range("H6").Value = -10.33984429491822

$ws.Range("G7").Value = 0.2406162607008728
$ws.Range("H7").Value = 46.87669836145837

$ws.Range("G8").Value = 0.07863868342346461
$ws.Range("H8").Value = -52.42925824282936

$ws.Range("G9").Value = 0.2005975577787466
$ws.Range("H9").Value = 2.821747521695436

$ws.Range("G10").Value = -0.1458836242684984
$ws.Range("H10").Value = -155.2807490100163

$ws.Range("G11").Value = -0.1274916552796612
$ws.Range("H11").Value = -7.344830708231812

$ws.Range("G12").Value = 0.1819993308417833
$ws.Range("H12").Value = 14.44390882666194

$ws.Range("G13").Value = 0.2401270006740394
$ws.Range("H13").Value = 16.75822257977357

$ws.Range("G14").Value = 0.2015489830203613
$ws.Range("H14").Value = 6.429587007161763

$ws.Range("G15").Value = 0.1805776681175887
$ws.Range("H15").Value = -27.74159526494

$ws.Range("G16").Value = -0.04026130921775196
$ws.Range("H16").Value = -210.3646840278084

$ws.Range("G17").Value = 0.02735981661359092
$ws.Range("H17").Value = -22.86619259309407

$ws.Range("G18").Value = 0.08810829926943783
$ws.Range("H18").Value = -49.16268408191976

$ws.Range("G19").Value = 0.1030911084170129
$ws.Range("H19").Value = -18.02091843152165

$ws.Range("G20").Value = 0.1187388106597019
$ws.Range("H20").Value = 3.564534264307088

$ws.Range("G21").Value = 0.1048707132341026
$ws.Range("H21").Value = 4.454123094461163

$ws.Range("G22").Value = 0.08880714469767767
$ws.Range("H22").Value = -5.722026245360742

$ws.Range("G23").Value = 0.1347052698572944
$ws.Range("H23").Value = 24.16377999973965

$ws.Range("G24").Value = -0.2139099522638939
$ws.Range("H24").Value = -71.64378931415649

$ws.Range("G25").Value = -0.1925158560902323
$ws.Range("H25").Value = 13.45816167014933

$ws.Range("G26").Value = 0.1725951200400995
$ws.Range("H26").Value = 8.562804555578293

$ws.Range("G27").Value = 0.1660360584265225
$ws.Range("H27").Value = -17.17509213928041

$ws.Range("G28").Value = 0.02394492607401189
$ws.Range("H28").Value = 398.176633264287

$ws.Range("G29").Value = 0.04444965162038098
$ws.Range("H29").Value = 189.0565774279435
